# Issue no. 70 @1418675072
#
# Mark todo items 1037 (row 10), 1038 (row 12) and 1039 (row 13) as closed:
# stamp the "closed" date (column I) and hide the now-completed rows, then
# move the active selection down to the next visible item (E20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$closedDate = 41988  # 2014-12-15

$rowsClosed = @(10, 12, 13)
foreach ($r in $rowsClosed) {
    $ws.Cells.Item($r, 9).Value = $closedDate
    $ws.Rows.Item($r).Hidden = $true
}

$ws.Range("E20").Select()
